$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1373.2759
$ws.Range("I15").Value = 1373.2759
$ws.Range("K15").Value = 4119.8277
$ws.Range("M15").Value = -3950.8277

$ws.Range("H62").Value = 2037.7858
$ws.Range("J62").Value = 2204.8572
$ws.Range("L62").Value = 2204.8572
$ws.Range("N62").Value = -3452.8572

$ws.Range("H65").Value = 2037.7858
$ws.Range("J65").Value = 2204.8572
$ws.Range("L65").Value = 11024.286
$ws.Range("N65").Value = -17264.286

$ws.Range("H86").Value = 43481388
$ws.Range("I86").Value = 58826830
$ws.Range("K86").Value = 58826830
$ws.Range("M86").Value = -58825707

$ws.Range("H88").Value = 3557.1538
$ws.Range("J88").Value = 4529.778
$ws.Range("L88").Value = 4529.778
$ws.Range("N88").Value = -5341.778

$ws.Range("H89").Value = 43481388
$ws.Range("I89").Value = 58826830
$ws.Range("K89").Value = 294134150
$ws.Range("M89").Value = -294128534

$ws.Range("H91").Value = 3557.1538
$ws.Range("J91").Value = 4529.778
$ws.Range("L91").Value = 4529.778
$ws.Range("N91").Value = -7337.778

$ws.Range("H98").Value = 3621.3142
$ws.Range("I98").Value = 2355.6667
$ws.Range("K98").Value = 2355.6667
$ws.Range("M98").Value = -857.6667000000002

$ws.Range("H107").Value = 16667054
$ws.Range("I107").Value = 458.33334
$ws.Range("K107").Value = 458.33334
$ws.Range("M107").Value = 1461.66666

$ws.Range("H118").Value = 719
$ws.Range("I118").Value = 562.8570999999999
$ws.Range("K118").Value = 1688.5713
$ws.Range("M118").Value = -31.57129999999984

$ws.Range("H122").Value = 3621.3142
$ws.Range("I122").Value = 2355.6667
$ws.Range("K122").Value = 7067.000100000001
$ws.Range("M122").Value = -4617.000100000001

$ws.Range("H132").Value = 41277.28
$ws.Range("I132").Value = 44692.695
$ws.Range("K132").Value = 134078.085
$ws.Range("M132").Value = -131548.085

$ws.Range("H137").Value = 2223175.8
$ws.Range("I137").Value = 955.8788
$ws.Range("J137").Value = 8334281
$ws.Range("K137").Value = 2867.6364
$ws.Range("L137").Value = 25002843
$ws.Range("M137").Value = -317.6363999999999
$ws.Range("N137").Value = -25007943

$ws.Range("H138").Value = 5169.5938
$ws.Range("I138").Value = 15503.6
$ws.Range("J138").Value = 3255.889
$ws.Range("K138").Value = 46510.8
$ws.Range("L138").Value = 9767.667000000001
$ws.Range("M138").Value = -41370.8
$ws.Range("N138").Value = -20047.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 176922.11
$ws.Range("I32").Value = 182629.67
$ws.Range("K32").Value = 182629.67
$ws.Range("M32").Value = -182342.67

$ws.Range("H74").Value = 339726.2
$ws.Range("I74").Value = 1574.762
$ws.Range("J74").Value = 931491.2
$ws.Range("K74").Value = 1574.762
$ws.Range("L74").Value = 931491.2
$ws.Range("M74").Value = -700.7619999999999
$ws.Range("N74").Value = -933239.2

$ws.Range("H77").Value = 339726.2
$ws.Range("I77").Value = 1574.762
$ws.Range("J77").Value = 931491.2
$ws.Range("K77").Value = 7873.809999999999
$ws.Range("L77").Value = 4657456
$ws.Range("M77").Value = -3505.809999999999
$ws.Range("N77").Value = -4666192

$ws.Range("H97").Value = 3923.1072
$ws.Range("I97").Value = 4510.5
$ws.Range("K97").Value = 4510.5
$ws.Range("M97").Value = -4014.5

$ws.Range("H110").Value = 1976.3334
$ws.Range("I110").Value = 1848.375
$ws.Range("K110").Value = 1848.375
$ws.Range("M110").Value = 196.625

$ws.Range("H132").Value = 2449.9614
$ws.Range("I132").Value = 2154.4358
$ws.Range("J132").Value = 3336.5386
$ws.Range("K132").Value = 6463.307400000001
$ws.Range("L132").Value = 10009.6158
$ws.Range("M132").Value = -3933.307400000001
$ws.Range("N132").Value = -15069.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2472.1538
$ws.Range("I86").Value = 1326.5555
$ws.Range("K86").Value = 1326.5555
$ws.Range("M86").Value = -203.5554999999999

$ws.Range("H89").Value = 2472.1538
$ws.Range("I89").Value = 1326.5555
$ws.Range("K89").Value = 6632.7775
$ws.Range("M89").Value = -1016.7775

$ws.Range("H107").Value = 8061.9116
$ws.Range("I107").Value = 9970
$ws.Range("K107").Value = 9970
$ws.Range("M107").Value = -8050

$ws.Range("H134").Value = 22502480
$ws.Range("I134").Value = 2294.5
$ws.Range("J134").Value = 64288540
$ws.Range("K134").Value = 6883.5
$ws.Range("L134").Value = 192865620
$ws.Range("M134").Value = -4348.5
$ws.Range("N134").Value = -192870690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2306.9666
$ws.Range("I31").Value = 1691.1034
$ws.Range("J31").Value = 2883.0967
$ws.Range("K31").Value = 1691.1034
$ws.Range("L31").Value = 2883.0967
$ws.Range("M31").Value = -1396.1034
$ws.Range("N31").Value = -3473.0967

$ws.Range("H34").Value = 2306.9666
$ws.Range("I34").Value = 1691.1034
$ws.Range("J34").Value = 2883.0967
$ws.Range("K34").Value = 1691.1034
$ws.Range("L34").Value = 2883.0967
$ws.Range("M34").Value = -1489.1034
$ws.Range("N34").Value = -3287.0967

$ws.Range("H107").Value = 1349.4667
$ws.Range("I107").Value = 1319.625
$ws.Range("J107").Value = 1468.8334
$ws.Range("K107").Value = 1319.625
$ws.Range("L107").Value = 1468.8334
$ws.Range("M107").Value = 600.375
$ws.Range("N107").Value = -5308.8334

$ws.Range("H134").Value = 2313.6365
$ws.Range("I134").Value = 1774.1428
$ws.Range("J134").Value = 3257.75
$ws.Range("K134").Value = 5322.428400000001
$ws.Range("L134").Value = 9773.25
$ws.Range("M134").Value = -2787.428400000001
$ws.Range("N134").Value = -14843.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 907.8
$ws.Range("J34").Value = 1281
$ws.Range("L34").Value = 3843
$ws.Range("N34").Value = -4011

$ws.Range("H106").Value = 18735
$ws.Range("J106").Value = 18735
$ws.Range("L106").Value = 56205
$ws.Range("N106").Value = -58097

$ws.Range("H107").Value = 125000800
$ws.Range("I107").Value = 586
$ws.Range("J107").Value = 166667540
$ws.Range("K107").Value = 1758
$ws.Range("L107").Value = 500002620
$ws.Range("M107").Value = 162
$ws.Range("N107").Value = -500006460

$ws.Range("H121").Value = 3298.5
$ws.Range("J121").Value = 6000
$ws.Range("L121").Value = 18000
$ws.Range("N121").Value = -20620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 408.2973
$ws.Range("I97").Value = 368.92593
$ws.Range("K97").Value = 368.92593
$ws.Range("M97").Value = 127.07407

$ws.Range("H107").Value = 46012.547
$ws.Range("I107").Value = 111311
$ws.Range("J107").Value = 805.9231
$ws.Range("K107").Value = 111311
$ws.Range("L107").Value = 805.9231
$ws.Range("M107").Value = -109391
$ws.Range("N107").Value = -4645.9231

$ws.Range("H113").Value = 2068.625
$ws.Range("I113").Value = 1924.75
$ws.Range("J113").Value = 2212.5
$ws.Range("K113").Value = 1924.75
$ws.Range("L113").Value = 2212.5
$ws.Range("M113").Value = 245.25
$ws.Range("N113").Value = -6552.5

$ws.Range("H122").Value = 1633.2954
$ws.Range("J122").Value = 1611.2
$ws.Range("L122").Value = 4833.6
$ws.Range("N122").Value = -9733.6

$ws.Range("H132").Value = 584693.75
$ws.Range("I132").Value = 2080.2856
$ws.Range("J132").Value = 856580
$ws.Range("K132").Value = 6240.8568
$ws.Range("L132").Value = 2569740
$ws.Range("M132").Value = -3710.8568
$ws.Range("N132").Value = -2574800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9810.862999999999
$ws.Range("I46").Value = 12203.333
$ws.Range("K46").Value = 12203.333
$ws.Range("M46").Value = -12015.333

$ws.Range("H55").Value = 1099.3617
$ws.Range("I55").Value = 1031.7241
$ws.Range("K55").Value = 1031.7241
$ws.Range("M55").Value = -858.7240999999999

$ws.Range("H76").Value = 11285.667
$ws.Range("I76").Value = 15000
$ws.Range("J76").Value = 9428.5
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 9428.5
$ws.Range("M76").Value = -14662
$ws.Range("N76").Value = -10104.5

$ws.Range("H79").Value = 11285.667
$ws.Range("I79").Value = 15000
$ws.Range("J79").Value = 9428.5
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 9428.5
$ws.Range("M79").Value = -13830
$ws.Range("N79").Value = -11768.5

$ws.Range("H82").Value = 1083.1765
$ws.Range("I82").Value = 1310.2858
$ws.Range("J82").Value = 924.2
$ws.Range("K82").Value = 1310.2858
$ws.Range("L82").Value = 924.2
$ws.Range("M82").Value = -949.2858000000001
$ws.Range("N82").Value = -1646.2

$ws.Range("H85").Value = 1083.1765
$ws.Range("I85").Value = 1310.2858
$ws.Range("J85").Value = 924.2
$ws.Range("K85").Value = 1310.2858
$ws.Range("L85").Value = 924.2
$ws.Range("M85").Value = -62.28580000000011
$ws.Range("N85").Value = -3420.2

$ws.Range("H93").Value = 1143.2174
$ws.Range("I93").Value = 1119.75
$ws.Range("K93").Value = 1119.75
$ws.Range("M93").Value = 128.25

$ws.Range("H132").Value = 3032.8286
$ws.Range("I132").Value = 2649.7222
$ws.Range("J132").Value = 3438.4707
$ws.Range("K132").Value = 7949.1666
$ws.Range("L132").Value = 10315.4121
$ws.Range("M132").Value = -5419.1666
$ws.Range("N132").Value = -15375.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H45").Value = 25956.2
$ws.Range("J45").Value = 27722.25
$ws.Range("L45").Value = 27722.25
$ws.Range("N45").Value = -28704.25

$ws.Range("H100").Value = 849.6667
$ws.Range("I100").Value = 843.375
$ws.Range("K100").Value = 1686.75
$ws.Range("M100").Value = -1145.75

$ws.Range("H136").Value = 45001.61
$ws.Range("I136").Value = 77505.08
$ws.Range("K136").Value = 232515.24
$ws.Range("M136").Value = -229965.24
